# Dane przykladowe dla pierwszych 20 atrybutow.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Arkusz1")
$ws2 = $wb.Worksheets.Item("Arkusz2")

# --- Fill sample data into Arkusz2 ---

$headers = @("IdPacjent","Wizyta ","Plec","Wiek","CzasTrwaniaChoroby","WiekZachorowania", `
    "ObjawyObecnieWyst","Drzenie","Sztywnosc","Spowolnienie","DyskinezyObecnie", `
    "DyskinezyOdLat","FluktuacjeObecnie","FluktuacjeOdLat","CzasOFF","CzasDyskinez", `
    "PoprawaPoLDopie","LDopaObecnie","AgonistaObecnie","LekiInne")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $ws2.Cells.Item(3, $col).Value = $headers[$i]
}

$data = @(
    @(1,1,1,34,2,31,2,0,1,0,1,1,1,1,2,1,1,10,$null,0),
    @(1,2,1,35,3,31,0,0,0,0,0,$null,0,$null,0,$null,1,10,$null,3),
    @(2,0,1,59,2,57,1,1,0,0,1,$null,1,0,1,2,0,15,20,4),
    @(2,1,0,60,3,57,1,1,0,0,1,1,1,1,1,2,1,15,20,4),
    @(1,1,0,70,1,69,4,0,0,0,0,$null,0,$null,2,$null,0,0,$null,5)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = 4 + $r
    $rowData = $data[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val) {
            $ws2.Cells.Item($row, 2 + $c).Value = $val
        }
    }
}

# Apply the "no border" bordered style (A1:A23 plus used range in the new rows)
# to column A of rows 1-23, matching the xf with applyBorder="1".
$ws2.Range("A1:U23").Borders.LineStyle = -4142

# Select the new data-entry cell used range
$ws2.Range("A11:D30").Select()
$ws1.Range("A11:D30").Select()

# Set selection on sheet2 and sheet1 per diff
$ws2.Range("J15").Select()

# Activate Arkusz2 as the active sheet/tab
$ws2.Activate()
$ws2.Range("J15").Select()

$ws1.Range("A11:D30").Select()
